# [Documentation] Weekly Report - 10
#
# "III. Next Week Plan" table: the Deadline column (D26:D31) is updated
# from the old date entry to the new deadline "13/11/2022" (entered as
# text, matching how it now appears in the report).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wx")

$ws.Range("D26:D31").Value = "13/11/2022"

# Leave the selection on the last edited cell.
$ws.Range("D31").Select()
